$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 142, shifting existing rows 142:200 down to 143:201
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new record
$ws.Range("A142").Value = 5
$ws.Range("B142").Value = "Macroferia Regional de Talca"
$ws.Range("C142").Value = "Maule"
$ws.Range("D142").Value = 44992
$ws.Range("E142").Value = 7
$ws.Range("F142").Value = 100112030
$ws.Range("G142").Value = "Poroto granado"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 200
$ws.Range("K142").Value = 30000
$ws.Range("L142").Value = 30000
$ws.Range("M142").Value = 30000
$ws.Range("N142").Value = "$/saco 25 kilos"
$ws.Range("O142").Value = "Región del Maule"
$ws.Range("P142").Value = 1200
$ws.Range("Q142").Value = 25
$ws.Range("R142").Value = "Hortaliza"
